# Update "想去人数" (F column) figures for both the 展览 sheet and the
# aggregated 全部类型 sheet, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1130
$ws1.Range("F4").Value = 255
$ws1.Range("F6").Value = 12110
$ws1.Range("F8").Value = 90
$ws1.Range("F9").Value = 11881
$ws1.Range("F10").Value = 4775
$ws1.Range("F11").Value = 585
$ws1.Range("F12").Value = 82
$ws1.Range("F13").Value = 28
$ws1.Range("F14").Value = 422
$ws1.Range("F16").Value = 935
$ws1.Range("F18").Value = 160

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1130
$ws4.Range("F4").Value = 255
$ws4.Range("F8").Value = 12110
$ws4.Range("F10").Value = 90
$ws4.Range("F11").Value = 11881
$ws4.Range("F12").Value = 4775
$ws4.Range("F13").Value = 585
$ws4.Range("F14").Value = 82
$ws4.Range("F15").Value = 28
$ws4.Range("F16").Value = 422
$ws4.Range("F18").Value = 935
$ws4.Range("F20").Value = 160

$wb.Save()
